$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status moved from "Ready for handoff" to "Handed back: in sync with en-US"
#    everywhere it appears (Overview summary columns + per-language Status col).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. The handback actually happened - stamp the real "Latest Handback
#    DateTime" (column H) in place of the "0001-01-01 00:00:00" placeholder.
#    zh-cn and de-de were handed back at different times.
# ---------------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-22 20:50:42"
$wsZhCn.Range("H3").Value = "2016-03-22 20:50:42"

$wsDeDe.Range("H2").Value = "2016-03-22 20:50:54"
$wsDeDe.Range("H3").Value = "2016-03-22 20:50:54"

# ---------------------------------------------------------------------------
# 3. Populate the new "Latest Target File" (F) and "Latest Handback File" (G)
#    columns for every row, each a hyperlink mirroring the corresponding
#    handoff file (F mirrors A, the source .md; G mirrors D, the handoff
#    .xlf) now that a target/handback round-trip exists.
# ---------------------------------------------------------------------------

# zh-cn sheet -------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/11e46888cb0954d172773ed23be311c2e0f63b66/e2e/0bb40168-d008-4845-bb69-20d061646237.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0bb40168-d008-4845-bb69-20d061646237.md") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/baf625206f21068775668097e622b6c03fd669b3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/11e46888cb0954d172773ed23be311c2e0f63b66/e2e/8a655b4d-52ca-4cc7-af93-86b65082ca1e.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/baf625206f21068775668097e622b6c03fd669b3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf") | Out-Null

# de-de sheet ---------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/11e46888cb0954d172773ed23be311c2e0f63b66/e2e/0bb40168-d008-4845-bb69-20d061646237.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0bb40168-d008-4845-bb69-20d061646237.md") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/707f60168b6f6b5cdae54a5f0691da7af0050c03/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/11e46888cb0954d172773ed23be311c2e0f63b66/e2e/8a655b4d-52ca-4cc7-af93-86b65082ca1e.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/707f60168b6f6b5cdae54a5f0691da7af0050c03/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf") | Out-Null
